$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Armada"
$ws.Range("B2").Value = "Fale"

$ws.Range("A3").Value = "Tony"
$ws.Range("B3").Value = "Ruso"

$ws.Range("A4").Value = "Palop"
$ws.Range("B4").Value = "Lope"

$ws.Range("A5").Value = "Kero"
$ws.Range("B5").Value = "Puche"

$ws.Range("A6").Value = "Kike"
$ws.Range("B6").Value = "Gonzo"

$ws.Range("A7").Value = "Papu"
$ws.Range("B7").Value = "Coquina"
